# The commit removed the "knowmads/ -> Page not found (404)" example row
# from the REST API specs sheet (the row's shared strings / values are
# cleared out, the row itself stays in place) and left the sheet's saved
# selection on the title row instead of cell A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 held: B10="knowmads/", C10="Page not found", E10=404, F10="N/A"
# (A10 and D10 were already blank). Clear the four populated cells so the
# whole row reads empty again, which also drops those two now-unused
# strings ("knowmads/", "Page not found") from the shared string table.
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

# Match the saved view state: selection moved from A11 to the header row.
$ws.Range("A1:F1").Select()
